$d = $word.ActiveDocument

$pairs = @(
    @("2024-01-16 Tuesday", "2024-01-17 Wednesday"),
    @("593×8=", "479×6="),
    @("844×4=", "279×3="),
    @("257×8=", "551×8="),
    @("804×9=", "951×5="),
    @("490×4=", "963×6="),
    @("975×6=", "707×2="),
    @("185×2=", "391×9="),
    @("568×3=", "363×9="),
    @("854×3=", "719×6="),
    @("972×8=", "352×7="),
    @("596×5=", "570×2="),
    @("121×9=", "681×8="),
    @("308×6=", "461×8="),
    @("843×3=", "320×7="),
    @("592×6=", "306×7="),
    @("141×7=", "596×8="),
    @("928×7=", "293×6="),
    @("707×9=", "999×5="),
    @("225×9=", "448×4="),
    @("935×4=", "655×2="),
    @("943×3=", "909×8="),
    @("506×5=", "194×6="),
    @("263×2=", "165×8="),
    @("635×2=", "469×3="),
    @("787×2=", "435×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
